$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-23T07:01:50.451543+00:00"
$ws.Range("K3").Value = "2025-11-23T07:01:50.451586+00:00"
$ws.Range("K4").Value = "2025-11-23T07:01:53.104468+00:00"
$ws.Range("K5").Value = "2025-11-23T07:01:53.104502+00:00"
$ws.Range("K6").Value = "2025-11-23T07:01:55.900574+00:00"
$ws.Range("K7").Value = "2025-11-23T07:01:58.212978+00:00"
$ws.Range("K8").Value = "2025-11-23T07:02:00.971759+00:00"
$ws.Range("K9").Value = "2025-11-23T07:02:00.971791+00:00"
$ws.Range("K10").Value = "2025-11-23T07:02:00.971811+00:00"
$ws.Range("K11").Value = "2025-11-23T07:02:03.701568+00:00"
$ws.Range("K12").Value = "2025-11-23T07:02:05.856455+00:00"
$ws.Range("K13").Value = "2025-11-23T07:02:08.128877+00:00"
$ws.Range("K14").Value = "2025-11-23T07:02:10.880350+00:00"
$ws.Range("K15").Value = "2025-11-23T07:02:13.665840+00:00"
$ws.Range("K16").Value = "2025-11-23T07:02:18.882252+00:00"
$ws.Range("K17").Value = "2025-11-23T07:02:18.882282+00:00"
$ws.Range("K18").Value = "2025-11-23T07:02:21.188261+00:00"
$ws.Range("K19").Value = "2025-11-23T07:02:21.188291+00:00"
$ws.Range("K20").Value = "2025-11-23T07:02:21.188309+00:00"
$ws.Range("K21").Value = "2025-11-23T07:02:24.030987+00:00"
$ws.Range("K22").Value = "2025-11-23T07:02:24.031019+00:00"
$ws.Range("K23").Value = "2025-11-23T07:02:26.756830+00:00"
$ws.Range("K24").Value = "2025-11-23T07:02:26.756861+00:00"
$ws.Range("K25").Value = "2025-11-23T07:02:26.756878+00:00"
$ws.Range("K26").Value = "2025-11-23T07:02:26.756896+00:00"
$ws.Range("K27").Value = "2025-11-23T07:02:29.129823+00:00"
$ws.Range("K28").Value = "2025-11-23T07:02:29.129852+00:00"
$ws.Range("K29").Value = "2025-11-23T07:02:31.580411+00:00"
$ws.Range("K30").Value = "2025-11-23T07:02:31.580440+00:00"
$ws.Range("K31").Value = "2025-11-23T07:02:31.580457+00:00"
$ws.Range("K32").Value = "2025-11-23T07:02:31.580472+00:00"
$ws.Range("K33").Value = "2025-11-23T07:02:33.874991+00:00"
$ws.Range("K34").Value = "2025-11-23T07:02:33.875021+00:00"
$ws.Range("K35").Value = "2025-11-23T07:02:38.606109+00:00"
$ws.Range("K36").Value = "2025-11-23T07:02:38.606140+00:00"
$ws.Range("K37").Value = "2025-11-23T07:02:41.377707+00:00"
$ws.Range("K38").Value = "2025-11-23T07:02:41.377735+00:00"
